$wb = $excel.ActiveWorkbook

# --- 1. Update status text: "Ready for handoff" -> "In Translation" ---
# This text appears in the Overview sheet (columns for each locale) as well
# as on each per-locale status sheet, so sweep every worksheet.
foreach ($ws in $wb.Worksheets) {
    [void]$ws.Cells.Replace("Ready for handoff", "In Translation")
}

# --- 2. Narrow the locale status columns from ~17.22 chars to ~13.41 chars ---
# The host's ColumnWidth setter rounds to whole pixels (2-decimal input,
# snapped to the sheet's Normal-font pixel grid), so the input below is the
# value that lands on the pixel closest to the target stored width
# (13.4101845877511 -> nearest achievable raw width 13.333333...).
$narrowWidth = 12.42

# Overview sheet: columns E (zh-cn) and F (de-de) get narrower.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $narrowWidth
$wsOverview.Columns.Item(6).ColumnWidth = $narrowWidth

# Per-locale sheets: column C ("Status") gets the same narrower width.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $narrowWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $narrowWidth
